# ---------------------------------------------------------------------------
# Reproduces the target diff:
#  - Insert a new leading "id" column (1..5) into "Лист1", shifting old data
#    right, and widen the header columns.
#  - Refresh the (cosmetic) selection on "Лист2".
#  - Insert a brand-new sheet "Лист8" (a single-cell signature/link sheet)
#    right before "Лист6", make it the active tab.
#  - Refresh selections on "Лист6", "Лист3", "Лист5" (pure cosmetic).
#  - Widen column C on "Лист4" and refresh its selection.
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# --- 1. "Лист1": insert an id column in front of the existing data --------
$ws1 = $wb.Worksheets.Item("Лист1")
$ws1.Columns("A").Insert()

$ws1.Range("A1").Value = "id"
$ws1.Range("A2").Value = 1
$ws1.Range("A3").Value = 2
$ws1.Range("A4").Value = 3
$ws1.Range("A5").Value = 4
$ws1.Range("A6").Value = 5

# narrow id column, widen the two text columns to match
$ws1.Columns("A").ColumnWidth = 1.83
$ws1.Range("B1:C6").Columns.ColumnWidth = 28.17

$ws1.Range("A2:C6").Select()

# --- 2. "Лист2": selection only -------------------------------------------
$ws2 = $wb.Worksheets.Item("Лист2")
$ws2.Range("A2:I6").Select()

# --- 3. Add the new "Лист8" sheet, positioned right before "Лист6" --------
# Worksheets.Add() without args drops the new sheet right before whatever is
# currently the ActiveSheet, so we first add a throwaway placeholder (this
# reserves the lower sheetId) and then add the real sheet in front of
# "Лист6" - this way the real sheet ends up with the next sheetId instead of
# the placeholder.
$placeholder = $wb.Worksheets.Add()
$placeholderName = $placeholder.Name

$newSheet = $wb.Worksheets.Add($wb.Worksheets.Item("Лист6"))
$newSheet.Name = "Лист8"

$wb.Worksheets.Item($placeholderName).Delete()

$wb.Worksheets.Item("Лист8").Range("A1").Value = "https://github.com/Ranimeboy/ra1"

# --- 4. "Лист6": selection only --------------------------------------------
$ws6 = $wb.Worksheets.Item("Лист6")
$ws6.Range("A2:B5").Select()

# --- 5. "Лист3": selection only --------------------------------------------
$ws3 = $wb.Worksheets.Item("Лист3")
$ws3.Range("A2:E17").Select()

# --- 6. "Лист4": widen column C, refresh selection --------------------------
$ws4 = $wb.Worksheets.Item("Лист4")
$ws4.Columns("C").ColumnWidth = 32.6
$ws4.Range("A2:C5").Select()

# --- 7. "Лист5": selection only ---------------------------------------------
$ws5 = $wb.Worksheets.Item("Лист5")
$ws5.Range("A2:E6").Select()

# --- 8. Finish on "Лист8" so it becomes the active tab (matches activeTab) --
$wb.Worksheets.Item("Лист8").Range("D4").Select()
